$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" header in column H, reusing the formatting of the other
# header cells (bold font, border, centered alignment) by copying G1's
# format into H1 before setting its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Add the numeric values for the new Save column
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
